$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the used range) for forcing text-typed values on numeric-looking
# strings, so Excel does not auto-convert them to floating point numbers. The helper
# cell is formatted as Text ("@"), the string is written there, copied, and pasted as
# values into the destination cell (pure value paste - no formatting is transferred).
$helper = $ws.Range('ZZ1')
$helper.NumberFormat = "@"

$helper.Value = '217.09'
$helper.Copy() | Out-Null
$ws.Range('D5').PasteSpecial(-4163) | Out-Null

$helper.Value = '19.04'
$helper.Copy() | Out-Null
$ws.Range('D10').PasteSpecial(-4163) | Out-Null

$helper.Value = '0.0843'
$helper.Copy() | Out-Null
$ws.Range('D11').PasteSpecial(-4163) | Out-Null

$helper.Value = '0.523'
$helper.Copy() | Out-Null
$ws.Range('D15').PasteSpecial(-4163) | Out-Null

$helper.Value = '64.32'
$helper.Copy() | Out-Null
$ws.Range('D16').PasteSpecial(-4163) | Out-Null

$helper.Value = '210.83'
$helper.Copy() | Out-Null
$ws.Range('D20').PasteSpecial(-4163) | Out-Null

$helper.Value = '6.17'
$helper.Copy() | Out-Null
$ws.Range('D22').PasteSpecial(-4163) | Out-Null

$helper.Value = '2.30'
$helper.Copy() | Out-Null
$ws.Range('D23').PasteSpecial(-4163) | Out-Null

$helper.Value = '145.94'
$helper.Copy() | Out-Null
$ws.Range('D25').PasteSpecial(-4163) | Out-Null

$helper.Value = '7.05'
$helper.Copy() | Out-Null
$ws.Range('D28').PasteSpecial(-4163) | Out-Null

$helper.Value = '15.51'
$helper.Copy() | Out-Null
$ws.Range('D29').PasteSpecial(-4163) | Out-Null

$helper.Value = '0.0503'
$helper.Copy() | Out-Null
$ws.Range('D30').PasteSpecial(-4163) | Out-Null

$helper.Value = '1.18'
$helper.Copy() | Out-Null
$ws.Range('D31').PasteSpecial(-4163) | Out-Null

$helper.Value = '3.35'
$helper.Copy() | Out-Null
$ws.Range('D32').PasteSpecial(-4163) | Out-Null

$helper.Value = '2.43'
$helper.Copy() | Out-Null
$ws.Range('D36').PasteSpecial(-4163) | Out-Null

$helper.Value = '0.526'
$helper.Copy() | Out-Null
$ws.Range('D38').PasteSpecial(-4163) | Out-Null

$helper.Value = '5.26'
$helper.Copy() | Out-Null
$ws.Range('D44').PasteSpecial(-4163) | Out-Null

$helper.Value = '91.21'
$helper.Copy() | Out-Null
$ws.Range('D45').PasteSpecial(-4163) | Out-Null

$helper.Value = '60.24'
$helper.Copy() | Out-Null
$ws.Range('D46').PasteSpecial(-4163) | Out-Null

$helper.Value = '1.58'
$helper.Copy() | Out-Null
$ws.Range('D47').PasteSpecial(-4163) | Out-Null

$helper.Clear() | Out-Null
$excel.CutCopyMode = 0

# Non-numeric-looking D values and all E (percentage) values can be set directly;
# Excel keeps these as plain text since they are not parseable as numbers.

$ws.Range('D2').Value = '26.705.10'
$ws.Range('E2').Value = '  -0.28%  '

$ws.Range('D3').Value = '1.635.63'
$ws.Range('E3').Value = '  -0.85%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('E5').Value = '  +0.28%  '

$ws.Range('E6').Value = '  -1.10%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  -0.94%  '

$ws.Range('E9').Value = '  -0.90%  '

$ws.Range('E10').Value = '  -0.93%  '

$ws.Range('E11').Value = '  +0.09%  '

$ws.Range('D12').Value = '1.863.54'
$ws.Range('E12').Value = '  -0.79%  '

$ws.Range('D13').Value = '1.630.19'
$ws.Range('E13').Value = '  -3.01%  '

$ws.Range('E14').Value = '  -1.30%  '

$ws.Range('E15').Value = '  -1.72%  '

$ws.Range('E16').Value = '  -1.86%  '

$ws.Range('D17').Value = '26.676.42'
$ws.Range('E17').Value = '  -0.45%  '

$ws.Range('E18').Value = '  -2.83%  '

$ws.Range('E19').Value = '  -0.02%  '

$ws.Range('E20').Value = '  -3.70%  '

$ws.Range('E21').Value = '  -1.21%  '

$ws.Range('E22').Value = '  -1.64%  '

$ws.Range('E23').Value = '  -2.42%  '

$ws.Range('E24').Value = '  -2.94%  '

$ws.Range('E25').Value = '  -0.36%  '

$ws.Range('E26').Value = '  -0.10%  '

$ws.Range('E27').Value = '  -2.20%  '

$ws.Range('E28').Value = '  -0.99%  '

$ws.Range('E29').Value = '  -1.55%  '

$ws.Range('E30').Value = '  -2.56%  '

$ws.Range('E31').Value = '  +0.45%  '

$ws.Range('E32').Value = '  -0.40%  '

$ws.Range('E33').Value = '  -1.73%  '

$ws.Range('D34').Value = '1.271.00'
$ws.Range('E34').Value = '  -0.89%  '

$ws.Range('E35').Value = '  -1.73%  '

$ws.Range('E36').Value = '  +0.10%  '

$ws.Range('E37').Value = '  -2.23%  '

$ws.Range('E38').Value = '  -1.86%  '

$ws.Range('E39').Value = '  -2.81%  '

$ws.Range('E40').Value = '  -0.10%  '

$ws.Range('E41').Value = '  -1.66%  '

$ws.Range('E42').Value = '  -2.62%  '

$ws.Range('D43').Value = '1.773.63'
$ws.Range('E43').Value = '  -0.81%  '

$ws.Range('E44').Value = '  -3.70%  '

$ws.Range('E45').Value = '  -0.91%  '

$ws.Range('E46').Value = '  +0.71%  '

$ws.Range('E47').Value = '  -2.19%  '

$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('E48').Value = '  -3.88%  '

$ws.Range('E49').Value = '  +0.45%  '

$ws.Range('E50').Value = '  -0.95%  '

$ws.Range('E51').Value = '  -3.00%  '
